# Generate Report for Handoff
# Refresh the "8d94f427-c817-4b74-878b-b6b6dfbbc76d" entity's latest handoff
# timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-17-18 04:17:44"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-18 04:17:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-18 04:17:44"
